$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.659.18"
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("D3").Value = "3.888.42"
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'604.11"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("D6").Value = "'170.46"
$ws.Range("E6").Value = "  +4.23%  "
$ws.Range("D7").Value = "3.888.61"
$ws.Range("E7").Value = "  +1.05%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +0.91%  "
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").Value = "'6.39"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("E12").Value = "  +1.49%  "
$ws.Range("E13").Value = "  +4.38%  "
$ws.Range("D14").Value = "'38.16"
$ws.Range("E14").Value = "  +3.39%  "
$ws.Range("D15").Value = "4.544.41"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "3.892.61"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "69.725.47"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "'18.71"
$ws.Range("E18").Value = "  +9.28%  "
$ws.Range("D19").Value = "'7.63"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("E20").Value = "  -0.90%  "
$ws.Range("D21").Value = "'11.11"
$ws.Range("E21").Value = "  -1.42%  "
$ws.Range("D22").Value = "'489.17"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").Value = "'0.746"
$ws.Range("E23").Value = "  +3.88%  "
$ws.Range("E24").Value = "  +1.69%  "
$ws.Range("D25").Value = "'85.22"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").Value = "'2.30"
$ws.Range("E26").Value = "  +2.75%  "
$ws.Range("D27").Value = "'12.33"
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").Value = "'10.11"
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("D31").Value = "4.039.24"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").Value = "'7.85"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("D35").Value = "3.857.38"
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").Value = "'6.11"
$ws.Range("E37").Value = "  +3.90%  "
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "'3.40"
$ws.Range("E38").Value = "  +14.30%  "
$ws.Range("E39").Value = "  +1.61%  "
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "'1.03"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  +2.90%  "
$ws.Range("E43").Value = "  +4.68%  "
$ws.Range("D44").Value = "'434.62"
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("D46").Value = "'8.67"
$ws.Range("E46").Value = "  +3.34%  "
$ws.Range("D48").Value = "'0.000277"
$ws.Range("E48").Value = "  +21.83%  "
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("D50").Value = "'143.60"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").Value = "'40.34"
$ws.Range("E51").Value = "  +4.06%  "
